$d = $word.ActiveDocument

# The template's ${imageN} placeholders (paragraphs 4-7, i.e. image2..image5)
# are being renumbered: 2->3, 3->5, 4->2, 5->4. Each replaced placeholder is
# rewritten as three separate runs - "${image", "<digit>", "}" - matching the
# target template's run layout. Paragraph indices are used (not text search)
# because the placeholder text itself is being shuffled/overlapping.
$w_ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Ordered list of (paragraph index, new digit) pairs, applied in document order.
$replacements = @(
    @{ Index = 4; Digit = "3" },
    @{ Index = 5; Digit = "5" },
    @{ Index = 6; Digit = "2" },
    @{ Index = 7; Digit = "4" }
)

foreach ($item in $replacements) {
    $p = $d.Paragraphs($item.Index)
    $rng = $p.Range
    $xml = '<w:p xmlns:w="' + $w_ns + '"><w:r><w:t>${image</w:t></w:r><w:r><w:t>' + $item.Digit + '</w:t></w:r><w:r><w:t>}</w:t></w:r></w:p>'
    [void]$rng.InsertXML($xml)
}
